$wb = $excel.ActiveWorkbook

# Sheet NextBus1 (index 1)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 45697.59561342592
$ws.Range("O2").Value = 13
$ws.Range("F3").Value = 45697.58855324074
$ws.Range("O3").Value = 3
$ws.Range("F4").Value = 45697.59538194445
$ws.Range("I4").Value = "SDA"
$ws.Range("O4").Value = 12
$ws.Range("F5").Value = 45697.59074074074
$ws.Range("O5").Value = 6
$ws.Range("F6").Value = 45697.58658564815
$ws.Range("O6").Value = 0
$ws.Range("F7").Value = 45697.5895949074
$ws.Range("O7").Value = 4
$ws.Range("F8").Value = 45697.5949537037
$ws.Range("O8").Value = 12
$ws.Range("F9").Value = 45697.59130787037
$ws.Range("O9").Value = 7
$ws.Range("F10").Value = 45697.58684027778
$ws.Range("O10").Value = 0
$ws.Range("F11").Value = 45697.58834490741
$ws.Range("O11").Value = 2
$ws.Range("F12").Value = 45697.58873842593
$ws.Range("O12").Value = 3
$ws.Range("F13").Value = 45697.59038194444
$ws.Range("F14").Value = 45697.59025462963
$ws.Range("O14").Value = 5
$ws.Range("F15").Value = 45697.59179398148
$ws.Range("O15").Value = 7

# Sheet NextBus2 (index 2)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 45697.60439814815
$ws.Range("O2").Value = 25
$ws.Range("F3").Value = 45697.59648148148
$ws.Range("O3").Value = 14
$ws.Range("F4").Value = 45697.60179398148
$ws.Range("O4").Value = 22
$ws.Range("F5").Value = 45697.59704861111
$ws.Range("O5").Value = 15
$ws.Range("F6").Value = 45697.59481481482
$ws.Range("L6").Value = "DD"
$ws.Range("O6").Value = 12
$ws.Range("F7").Value = 45697.59584490741
$ws.Range("I7").Value = "SEA"
$ws.Range("F8").Value = 45697.60163194445
$ws.Range("O8").Value = 21
$ws.Range("F9").Value = 45697.59827546297
$ws.Range("O9").Value = 17
$ws.Range("F10").Value = 45697.59410879629
$ws.Range("O10").Value = 11
$ws.Range("F11").Value = 45697.599375
$ws.Range("L11").Value = "DD"
$ws.Range("O11").Value = 18
$ws.Range("F12").Value = 45697.59513888889
$ws.Range("O12").Value = 12
$ws.Range("F13").Value = 45697.5969675926
$ws.Range("L13").Value = "SD"
$ws.Range("O13").Value = 15
$ws.Range("F14").Value = 45697.59560185186
$ws.Range("F15").Value = 45697.59932870371
$ws.Range("O15").Value = 18

# Sheet NextBus3 (index 3)
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 45697.61236111111
$ws.Range("J2").Value = 0
$ws.Range("O2").Value = 37
$ws.Range("F3").Value = 45697.60311342592
$ws.Range("O3").Value = 24
$ws.Range("F4").Value = 45697.60972222222
$ws.Range("O4").Value = 33
$ws.Range("F5").Value = 45697.60918981482
$ws.Range("O5").Value = 32
$ws.Range("F6").Value = 45697.59932870371
$ws.Range("O6").Value = 18
$ws.Range("F7").Value = 45697.60267361111
$ws.Range("O7").Value = 23
$ws.Range("F8").Value = 45697.61180555556
$ws.Range("O8").Value = 36
$ws.Range("F9").Value = 45697.61460648148
$ws.Range("O9").Value = 40
$ws.Range("F10").Value = 45697.61116898148
$ws.Range("O10").Value = 35
$ws.Range("F11").Value = 45697.60493055556
$ws.Range("L11").Value = "DD"
$ws.Range("O11").Value = 26
$ws.Range("F12").Value = 45697.60075231481
$ws.Range("L12").Value = "SD"
$ws.Range("O12").Value = 20
$ws.Range("F13").Value = 45697.60748842593
$ws.Range("L13").Value = "SD"
$ws.Range("O13").Value = 30
$ws.Range("F14").Value = 45697.60258101852
$ws.Range("O14").Value = 23
$ws.Range("F15").Value = 45697.61060185185
$ws.Range("J15").Value = 0
$ws.Range("O15").Value = 34
